$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (375-385), continuing the daily series through 2021-09-20,
# as per commit "aggiornamento fino a 20/09/2021".
$rows = @(
    @{ Row = 375; Date = 44449; B = 1; C = 8;  D = 51.66957308015243 },
    @{ Row = 376; Date = 44450; B = 4; C = 8;  D = 51.66957308015243 },
    @{ Row = 377; Date = 44451; B = 5; C = 12; D = 77.50435962022863 },
    @{ Row = 378; Date = 44452; B = 0; C = 11; D = 71.04566298520959 },
    @{ Row = 379; Date = 44453; B = 0; C = 10; D = 64.58696635019054 },
    @{ Row = 380; Date = 44454; B = 0; C = 10; D = 64.58696635019054 },
    @{ Row = 381; Date = 44455; B = 0; C = 10; D = 64.58696635019054 },
    @{ Row = 382; Date = 44456; B = 1; C = 10; D = 64.58696635019054 },
    @{ Row = 383; Date = 44457; B = 0; C = 6;  D = 38.75217981011431 },
    @{ Row = 384; Date = 44458; B = 1; C = 2;  D = 12.91739327003811 },
    @{ Row = 385; Date = 44459; B = 3; C = 5;  D = 32.29348317509527 }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Column A: date serial, with the same style (s="2", date number format)
    # used by the cell directly above it. Copy/PasteSpecial clones the
    # existing style exactly instead of minting a brand-new cellXf.
    $ws.Cells.Item($row - 1, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws.Cells.Item($row, 1).Value = $r.Date

    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
}
